# Update data: 3 December 2021
# Adds the 2021-11-01 (serial 44501) observation to both the "Canada"
# sheet (sheet1) and the "Province" sheet (sheet2).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Canada")
$ws2 = $wb.Worksheets.Item("Province")

# ---------------------------------------------------------------------
# Sheet "Canada": new row 24
# ---------------------------------------------------------------------
$ws1.Cells.Item(24, 1).Value = 44501
$ws1.Range("A24").NumberFormat = "d-mmm-yy"
$ws1.Cells.Item(24, 2).Value = "Canada"
$ws1.Range("B24").NumberFormat = "d-mmm-yy"
$ws1.Cells.Item(24, 4).Value = 1243.8
$ws1.Cells.Item(24, 5).Value = 1195.3
$ws1.Cells.Item(24, 3).Formula = "=(D24-E24)/E24*100"

# ---------------------------------------------------------------------
# Sheet "Province": new rows 222-231
# ---------------------------------------------------------------------
$ws2.Cells.Item(222, 1).Value = 44501
$ws2.Range("A222").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(222, 2).Value = "Newfoundland & Labrador"
$ws2.Range("B222").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(222, 4).Value = 26.6
$ws2.Cells.Item(222, 5).Value = 29.8
$ws2.Cells.Item(222, 3).Formula = "=(D222-E222)/E222*100"

$ws2.Cells.Item(223, 1).Value = 44501
$ws2.Range("A223").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(223, 2).Value = "Prince Edward Island"
$ws2.Cells.Item(223, 4).Value = 7.2
$ws2.Cells.Item(223, 5).Value = 6.8
$ws2.Cells.Item(223, 3).Formula = "=(D223-E223)/E223*100"

$ws2.Cells.Item(224, 1).Value = 44501
$ws2.Range("A224").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(224, 2).Value = "Nova Scotia"
$ws2.Cells.Item(224, 4).Value = 41.6
$ws2.Cells.Item(224, 5).Value = 40.1
$ws2.Cells.Item(224, 3).Formula = "=(D224-E224)/E224*100"

$ws2.Cells.Item(225, 1).Value = 44501
$ws2.Range("A225").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(225, 2).Value = "New Brunswick"
$ws2.Cells.Item(225, 4).Value = 34.1
$ws2.Cells.Item(225, 5).Value = 31.1
$ws2.Cells.Item(225, 3).Formula = "=(D225-E225)/E225*100"

$ws2.Cells.Item(226, 1).Value = 44501
$ws2.Range("A226").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(226, 2).Value = "Quebec"
$ws2.Cells.Item(226, 4).Value = 204.4
$ws2.Cells.Item(226, 5).Value = 253.9
$ws2.Cells.Item(226, 3).Formula = "=(D226-E226)/E226*100"

$ws2.Cells.Item(227, 1).Value = 44501
$ws2.Range("A227").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(227, 2).Value = "Ontario"
$ws2.Cells.Item(227, 4).Value = 519.1
$ws2.Cells.Item(227, 5).Value = 434.4
$ws2.Cells.Item(227, 3).Formula = "=(D227-E227)/E227*100"

$ws2.Cells.Item(228, 1).Value = 44501
$ws2.Range("A228").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(228, 2).Value = "Manitoba"
$ws2.Cells.Item(228, 4).Value = 35.4
$ws2.Cells.Item(228, 5).Value = 38.9
$ws2.Cells.Item(228, 3).Formula = "=(D228-E228)/E228*100"

$ws2.Cells.Item(229, 1).Value = 44501
$ws2.Range("A229").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(229, 2).Value = "Saskatchewan"
$ws2.Cells.Item(229, 4).Value = 30.6
$ws2.Cells.Item(229, 5).Value = 36.6
$ws2.Cells.Item(229, 3).Formula = "=(D229-E229)/E229*100"

$ws2.Cells.Item(230, 1).Value = 44501
$ws2.Range("A230").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(230, 2).Value = "Alberta"
$ws2.Cells.Item(230, 4).Value = 186.2
$ws2.Cells.Item(230, 5).Value = 182.3
$ws2.Cells.Item(230, 3).Formula = "=(D230-E230)/E230*100"

$ws2.Cells.Item(231, 1).Value = 44501
$ws2.Range("A231").NumberFormat = "d-mmm-yy"
$ws2.Cells.Item(231, 2).Value = "British Columbia"
$ws2.Cells.Item(231, 4).Value = 158.5
$ws2.Cells.Item(231, 5).Value = 141.19999999999999
$ws2.Cells.Item(231, 3).Formula = "=(D231-E231)/E231*100"

# ---------------------------------------------------------------------
# View state: move the active cell to the new last row on each sheet,
# and restore "Province" as the active/selected sheet (it was already
# the active tab before the edit).
# ---------------------------------------------------------------------
[void]$ws1.Range("A24").Select()
[void]$ws2.Range("D232").Select()
